$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data rows appended at the bottom of the "Diaria" table (rows 83-84),
# mirroring the existing rows above them (plain values, column A holding
# the date label as text, same as the rest of the column).

# Temporarily force column A to text formatting so the dd-mm-yyyy looking
# labels are stored as text (shared strings) instead of being auto-parsed
# into date serials by Excel's smart input - then drop the formatting
# again so the cells end up unstyled, exactly like the rest of column A.
$ws.Range("A83:A84").NumberFormat = "@"

# Row 83
$ws.Range("A83").Value = "02-11-2021"
$ws.Range("B83").Value = 50000
$ws.Range("C83").Value = 165000
$ws.Range("D83").Value = 50000
$ws.Range("E83").Value = 30000
$ws.Range("F83").Value = 20000
$ws.Range("G83").Value = 3.05

# Row 84
$ws.Range("A84").Value = "03-11-2021"
$ws.Range("B84").Value = 50000
$ws.Range("C84").Value = 102000
$ws.Range("D84").Value = 50000
$ws.Range("E84").Value = 33000
$ws.Range("F84").Value = 17000
$ws.Range("G84").Value = 3.1

# Drop the temporary text formatting again so A83:A84 stay unstyled.
$ws.Range("A83:A84").ClearFormats()
